$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 521.25
$ws.Range("I28").Value = 395
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 395
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = 90
$ws.Range("N28").Value = -1870
# Row 92
$ws.Range("H92").Value = 829.8125
$ws.Range("I92").Value = 816.7692
$ws.Range("J92").Value = 886.3333
$ws.Range("K92").Value = 816.7692
$ws.Range("L92").Value = 886.3333
$ws.Range("M92").Value = 431.2308
$ws.Range("N92").Value = -3382.3333
# Row 100
$ws.Range("H100").Value = 3278.889
$ws.Range("J100").Value = 900
$ws.Range("L100").Value = 900
$ws.Range("N100").Value = -1982
# Row 118
$ws.Range("H118").Value = 378.16666
$ws.Range("I118").Value = 378.16666
$ws.Range("K118").Value = 1134.49998
$ws.Range("M118").Value = 522.5000199999999
# Row 125
$ws.Range("H125").Value = 500
$ws.Range("I125").Value = 500
$ws.Range("K125").Value = 4500
$ws.Range("M125").Value = -2040
# Row 127
$ws.Range("H127").Value = 1426.75
$ws.Range("I127").Value = 995
$ws.Range("K127").Value = 2985
$ws.Range("M127").Value = 1975
# Row 131
$ws.Range("H131").Value = 4632.5
$ws.Range("J131").Value = 4600
$ws.Range("L131").Value = 13800
$ws.Range("N131").Value = -23880
# Row 132
$ws.Range("H132").Value = 1368.091
$ws.Range("I132").Value = 1368.091
$ws.Range("K132").Value = 4104.272999999999
$ws.Range("M132").Value = -1574.272999999999
# Row 138
$ws.Range("H138").Value = 2336.516
$ws.Range("I138").Value = 1338.8
$ws.Range("K138").Value = 4016.4
$ws.Range("M138").Value = 1123.6

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8258.857
$ws.Range("I32").Value = 8258.857
$ws.Range("K32").Value = 8258.857
$ws.Range("M32").Value = -7971.857
# Row 45
$ws.Range("H45").Value = 2916
$ws.Range("I45").Value = 2916
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2916
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2539
$ws.Range("N45").ClearContents()
# Row 88
$ws.Range("H88").Value = 4624.75
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 5166.3335
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 5166.3335
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -5978.3335
# Row 91
$ws.Range("H91").Value = 4624.75
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 5166.3335
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 5166.3335
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -7974.3335
# Row 97
$ws.Range("H97").Value = 261.75
$ws.Range("I97").Value = 261.75
$ws.Range("K97").Value = 261.75
$ws.Range("M97").Value = 234.25
# Row 110
$ws.Range("H110").Value = 3886.9167
$ws.Range("I110").Value = 2199.889
$ws.Range("J110").Value = 8948
$ws.Range("K110").Value = 2199.889
$ws.Range("L110").Value = 8948
$ws.Range("M110").Value = -154.8890000000001
$ws.Range("N110").Value = -13038
# Row 132
$ws.Range("H132").Value = 3200.353
$ws.Range("J132").Value = 3287
$ws.Range("L132").Value = 9861
$ws.Range("N132").Value = -14921

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 86
$ws.Range("H86").Value = 11666.667
# Row 89
$ws.Range("H89").Value = 11666.667
# Row 94
$ws.Range("H94").Value = 4024.7856
$ws.Range("I94").Value = 4141
$ws.Range("J94").Value = 3734.25
$ws.Range("K94").Value = 4141
$ws.Range("L94").Value = 3734.25
$ws.Range("M94").Value = -3690
$ws.Range("N94").Value = -4636.25
# Row 99
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 498
$ws.Range("N99").ClearContents()
# Row 102
$ws.Range("H102").Value = 35999.5
$ws.Range("J102").Value = 42000
$ws.Range("L102").Value = 42000
$ws.Range("N102").Value = -48490
# Row 103
$ws.Range("H103").Value = 120999
$ws.Range("J103").Value = 120999
$ws.Range("L103").Value = 120999
$ws.Range("N103").Value = -123343
# Row 105
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# Row 106
$ws.Range("H106").Value = 47557
$ws.Range("J106").Value = 47557
$ws.Range("L106").Value = 47557
$ws.Range("N106").Value = -50081
# Row 107
$ws.Range("H107").Value = 974.25
$ws.Range("I107").Value = 974.25
$ws.Range("K107").Value = 974.25
$ws.Range("M107").Value = 945.75

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 134
$ws.Range("H134").Value = 9805.272000000001
$ws.Range("I134").Value = 9991.375
$ws.Range("K134").Value = 29974.125
$ws.Range("M134").Value = -27439.125

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 309.44446
$ws.Range("I38").Value = 61.333332
$ws.Range("J38").Value = 433.5
$ws.Range("K38").Value = 183.999996
$ws.Range("L38").Value = 1300.5
$ws.Range("M38").Value = 163.000004
$ws.Range("N38").Value = -1994.5
# Row 126
$ws.Range("H126").Value = 6833.3335
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -5560
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 7333.3335
$ws.Range("I126").Value = 7333.3335
$ws.Range("K126").Value = 22000.0005
$ws.Range("M126").Value = -19530.0005

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2609.875
$ws.Range("I7").Value = 2554.4285
$ws.Range("J7").Value = 2998
$ws.Range("K7").Value = 2554.4285
$ws.Range("L7").Value = 2998
$ws.Range("M7").Value = -2442.4285
$ws.Range("N7").Value = -3222
# Row 122
$ws.Range("H122").Value = 9625
$ws.Range("I122").Value = 9625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 28875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -26425
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 2609.875
$ws.Range("I126").Value = 2554.4285
$ws.Range("J126").Value = 2998
$ws.Range("K126").Value = 7663.2855
$ws.Range("L126").Value = 8994
$ws.Range("M126").Value = -5193.2855
$ws.Range("N126").Value = -13934
# Row 132
$ws.Range("H132").Value = 4915.75
$ws.Range("I132").Value = 4166.1665
$ws.Range("J132").Value = 5665.3335
$ws.Range("K132").Value = 12498.4995
$ws.Range("L132").Value = 16996.0005
$ws.Range("M132").Value = -9968.499500000002
$ws.Range("N132").Value = -22056.0005
# Row 136
$ws.Range("H136").Value = 3578.5
$ws.Range("I136").Value = 3578.5
$ws.Range("K136").Value = 10735.5
$ws.Range("M136").Value = -8185.5

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 4669.2856
$ws.Range("I96").Value = 4669.2856
$ws.Range("K96").Value = 4669.2856
$ws.Range("M96").Value = -3296.2856
# Row 132
$ws.Range("H132").Value = 2473.7
$ws.Range("I132").Value = 1373.5
$ws.Range("J132").Value = 4124
$ws.Range("K132").Value = 4120.5
$ws.Range("L132").Value = 12372
$ws.Range("M132").Value = -1590.5
$ws.Range("N132").Value = -17432
# Row 136
$ws.Range("H136").Value = 6212.778
$ws.Range("I136").Value = 6219.3335
$ws.Range("K136").Value = 18658.0005
$ws.Range("M136").Value = -16108.0005
